# Auto-generated edit script: refreshes the cryptos.xlsx price/volume table
# per commit 'Updated cryptos list on Tue Feb 21 18:44:48 UTC 2023 with GitHub Actions'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 19 & 20 swap coin identity (Avalanche <-> TRON) in this refresh, in
# addition to getting new price/volume figures (set further below).
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"

# Price column (D) updates. These cells hold plain text (e.g. '1.001',
# '0.00001321'); a handful of the new values parse as valid numbers, so
# Excel would silently convert the cell to a Number (dropping trailing
# zeros / switching to scientific notation) unless we briefly mark the
# cell as Text first. The NumberFormat is restored to the sheet's normal
# (General) style right after so no visible formatting changes stick.
$ws.Range("D2").Value = "24.590.76"
$ws.Range("D3").Value = "1.676.20"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3927"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3957"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.002"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.404"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08666"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.342"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001321"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.726"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "1.670.63"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07024"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "21.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.093"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "24.602.59"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.367"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.785"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.868"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "160.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "146.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.315"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.510"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "1.840.45"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.03095"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08318"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.971"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2812"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9888"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09537"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.512"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7917"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7129"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.564"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.167"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08638"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.332"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.71"
$ws.Range("D51").Style = "Normal"

# Volume(1h) column (E) updates -- these are padded percentage strings
# (e.g. '  -1.19%  ') and are never number-like, so no text-forcing is
# required here.
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  -2.29%  "
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  -4.16%  "
$ws.Range("E11").Value = "  -5.62%  "
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("E16").Value = "  -4.01%  "
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("E18").Value = "  -2.76%  "
$ws.Range("E19").Value = "  -2.41%  "
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  -3.97%  "
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E26").Value = "  -3.67%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  -11.72%  "
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("E32").Value = "  +10.23%  "
$ws.Range("E33").Value = "  -6.09%  "
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("E35").Value = "  -4.82%  "
$ws.Range("E36").Value = "  -5.88%  "
$ws.Range("E37").Value = "  -1.69%  "
$ws.Range("E38").Value = "  -4.72%  "
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("E41").Value = "  -4.73%  "
$ws.Range("E42").Value = "  -6.92%  "
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("E44").Value = "  -6.40%  "
$ws.Range("E45").Value = "  -4.50%  "
$ws.Range("E46").Value = "  -5.34%  "
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("E48").Value = "  +3.42%  "
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("E50").Value = "  -5.26%  "
$ws.Range("E51").Value = "  -2.50%  "
